$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.8163891899954194
$ws.Range("C2").Value = 0.1363548090420679
$ws.Range("D2").Value = 0.112408321696833
$ws.Range("E2").Value = 0.135349934076336
$ws.Range("F2").Value = 2.131864999578376
$ws.Range("J2").Value = 0.190852553655283
$ws.Range("K2").Value = 0.4121628667787718
$ws.Range("L2").Value = 0.2484800205762738
$ws.Range("O2").Value = 5.740659178069734

# Row 3
$ws.Range("B3").Value = 0.7785231929974827
$ws.Range("C3").Value = 0.1356897337922049
$ws.Range("D3").Value = 0.1105904294790747
$ws.Range("E3").Value = 0.135445654879554
$ws.Range("F3").Value = 2.141136752964265
$ws.Range("J3").Value = 0.1919772027742077
$ws.Range("K3").Value = 0.3778670620411333
$ws.Range("L3").Value = 0.2445345036907725
$ws.Range("O3").Value = 5.773909290133986

# Row 4
$ws.Range("B4").Value = 0.7555567052488641
$ws.Range("C4").Value = 0.1352854131625705
$ws.Range("D4").Value = 0.1095130407305192
$ws.Range("E4").Value = 0.1355497092653319
$ws.Range("F4").Value = 2.14773974851235
$ws.Range("J4").Value = 0.1927271816332912
$ws.Range("K4").Value = 0.3568920994118514
$ws.Range("L4").Value = 0.2422055216403365
$ws.Range("O4").Value = 5.7968059251113

# Row 5
$ws.Range("B5").Value = 0.746269661794031
$ws.Range("C5").Value = 0.1351216867546796
$ws.Range("D5").Value = 0.1090838145701483
$ws.Range("E5").Value = 0.1356035281035144
$ws.Range("F5").Value = 2.150659589764643
$ws.Range("J5").Value = 0.1930477720003552
$ws.Range("K5").Value = 0.3483660032231199
$ws.Range("L5").Value = 0.241280089020151
$ws.Range("O5").Value = 5.806760723277506

# Row 6
$ws.Range("B6").Value = 0.7447319245638653
$ws.Range("C6").Value = 0.1350945634270815
$ws.Range("D6").Value = 0.1090131366983513
$ws.Range("E6").Value = 0.1356131550351307
$ws.Range("F6").Value = 2.151158268573923
$ws.Range("J6").Value = 0.1931019104364839
$ws.Range("K6").Value = 0.3469515584799154
$ws.Range("L6").Value = 0.2411278532071819
$ws.Range("O6").Value = 5.808451423383886

# Row 7
$ws.Range("B7").Value = 0.7554311644953771
$ws.Range("C7").Value = 0.1352832008601794
$ws.Range("D7").Value = 0.1095072122000573
$ws.Range("E7").Value = 0.1355503888234839
$ws.Range("F7").Value = 2.147778198746117
$ws.Range("J7").Value = 0.1927314445945427
$ws.Range("K7").Value = 0.3567770262327201
$ws.Range("L7").Value = 0.2421929450360736
$ws.Range("O7").Value = 5.79693765100086

# Row 8
$ws.Range("B8").Value = 0.8032745992379091
$ws.Range("C8").Value = 0.1361246659009652
$ws.Range("D8").Value = 0.1117734942705084
$ws.Range("E8").Value = 0.1353735578866022
$ws.Range("F8").Value = 2.134873135262637
$ws.Range("J8").Value = 0.1912280087187845
$ws.Range("K8").Value = 0.4003208068206163
$ws.Range("L8").Value = 0.2471002505861861
$ws.Range("O8").Value = 5.751609125904878

# Row 9
$ws.Range("B9").Value = 0.8993173966393044
$ws.Range("C9").Value = 0.1378059466874006
$ws.Range("D9").Value = 0.116523202973724
$ws.Range("E9").Value = 0.1353848802571562
$ws.Range("F9").Value = 2.116779194887343
$ws.Range("J9").Value = 0.1887504721535933
$ws.Range("K9").Value = 0.4863472631278682
$ws.Range("L9").Value = 0.2574616504627016
$ws.Range("O9").Value = 5.682390175139574

# Row 10
$ws.Range("B10").Value = 0.9712067489253684
$ws.Range("C10").Value = 0.1390591594661572
$ws.Range("D10").Value = 0.1201963663013714
$ws.Range("E10").Value = 0.1356100221079721
$ws.Range("F10").Value = 2.10787296037374
$ws.Range("J10").Value = 0.1872159293009439
$ws.Range("K10").Value = 0.5499193840906003
$ws.Range("L10").Value = 0.2655193970173428
$ws.Range("O10").Value = 5.643509214705944

# Row 11
$ws.Range("B11").Value = 1.004193762961307
$ws.Range("C11").Value = 0.1396329762818738
$ws.Range("D11").Value = 0.1219066955649595
$ws.Range("E11").Value = 0.1357592229804006
$ws.Range("F11").Value = 2.104772013469812
$ws.Range("J11").Value = 0.1865796027800393
$ws.Range("K11").Value = 0.5789161876818127
$ws.Range("L11").Value = 0.2692808023836335
$ws.Range("O11").Value = 5.628417951334768

# Row 12
$ws.Range("B12").Value = 1.016725329855916
$ws.Range("C12").Value = 0.1398507815740899
$ws.Range("D12").Value = 0.1225599617790323
$ws.Range("E12").Value = 0.1358224212378296
$ws.Range("F12").Value = 2.103734270774737
$ws.Range("J12").Value = 0.1863475007571083
$ws.Range("K12").Value = 0.5899072071808291
$ws.Range("L12").Value = 0.2707188374277649
$ws.Range("O12").Value = 5.623076264522922

# Row 13
$ws.Range("B13").Value = 1.014024661380859
$ws.Range("C13").Value = 0.1398038508273984
$ws.Range("D13").Value = 0.1224190210695895
$ws.Range("E13").Value = 0.1358085128649442
$ws.Range("F13").Value = 2.103951698212967
$ws.Range("J13").Value = 0.1863970942669511
$ws.Range("K13").Value = 0.5875396342049726
$ws.Range("L13").Value = 0.2704085246296586
$ws.Range("O13").Value = 5.624210104757111

# Row 14
$ws.Range("B14").Value = 1.005223943210098
$ws.Range("C14").Value = 0.1396508850919531
$ws.Range("D14").Value = 0.1219603282611814
$ws.Range("E14").Value = 0.1357642882397805
$ws.Range("F14").Value = 2.104683902682766
$ws.Range("J14").Value = 0.1865603301207734
$ws.Range("K14").Value = 0.5798202172127844
$ws.Range("L14").Value = 0.2693988370463245
$ws.Range("O14").Value = 5.627971012484721

# Row 15
$ws.Range("B15").Value = 0.9998384488266367
$ws.Range("C15").Value = 0.1395572553574027
$ws.Range("D15").Value = 0.1216800933215865
$ws.Range("E15").Value = 0.1357380709685856
$ws.Range("F15").Value = 2.10515017281368
$ws.Range("J15").Value = 0.1866614702342417
$ws.Range("K15").Value = 0.575093209099748
$ws.Range("L15").Value = 0.2687821513302993
$ws.Range("O15").Value = 5.630323251768317

# Row 16
$ws.Range("B16").Value = 0.9690566230878517
$ws.Range("C16").Value = 0.1390217322229859
$ws.Range("D16").Value = 0.1200853795254062
$ws.Range("E16").Value = 0.135601210458276
$ws.Range("F16").Value = 2.108094727519941
$ws.Range("J16").Value = 0.187258755673426
$ws.Range("K16").Value = 0.5480258844236801
$ws.Range("L16").Value = 0.2652754999739955
$ws.Range("O16").Value = 5.644547679711877

# Row 17
$ws.Range("B17").Value = 0.9502452119734812
$ws.Range("C17").Value = 0.1386941440702216
$ws.Range("D17").Value = 0.1191171182201742
$ws.Range("E17").Value = 0.1355292130799448
$ws.Range("F17").Value = 2.110144453695739
$ws.Range("J17").Value = 0.1876409725559967
$ws.Range("K17").Value = 0.5314403941826811
$ws.Range("L17").Value = 0.2631487607042544
$ws.Range("O17").Value = 5.653938603798593

# Row 18
$ws.Range("B18").Value = 0.9394521804697717
$ws.Range("C18").Value = 0.1385060758654006
$ws.Range("D18").Value = 0.118563911114876
$ws.Range("E18").Value = 0.1354922071839724
$ws.Range("F18").Value = 2.111412881729692
$ws.Range("J18").Value = 0.1878666263610178
$ws.Range("K18").Value = 0.521908191134969
$ws.Range("L18").Value = 0.2619345497112988
$ws.Range("O18").Value = 5.659584365652819

# Row 19
$ws.Range("B19").Value = 0.9358024708085964
$ws.Range("C19").Value = 0.1384424603097898
$ws.Range("D19").Value = 0.1183772439974149
$ws.Range("E19").Value = 0.1354804350613712
$ws.Range("F19").Value = 2.111857723036721
$ws.Range("J19").Value = 0.1879440277535203
$ws.Range("K19").Value = 0.5186820282753217
$ws.Range("L19").Value = 0.2615249937677504
$ws.Range("O19").Value = 5.661537898328845

# Row 20
$ws.Range("B20").Value = 0.9522449512086268
$ws.Range("C20").Value = 0.1387289801362996
$ws.Range("D20").Value = 0.1192198076485766
$ws.Range("E20").Value = 0.1355364216178465
$ws.Range("F20").Value = 2.109916997244937
$ws.Range("J20").Value = 0.1875996834464893
$ws.Range("K20").Value = 0.5332051933442301
$ws.Range("L20").Value = 0.2633742216802801
$ws.Range("O20").Value = 5.652913636230721

# Row 21
$ws.Range("B21").Value = 1.007807845237608
$ws.Range("C21").Value = 0.1396958010813094
$ws.Range("D21").Value = 0.1220949059284493
$ws.Range("E21").Value = 0.1357770964961773
$ws.Range("F21").Value = 2.104465132748132
$ws.Range("J21").Value = 0.1865121434495585
$ws.Range("K21").Value = 0.582087315211993
$ws.Range("L21").Value = 0.2696950366310062
$ws.Range("O21").Value = 5.626856220217491

# Row 22
$ws.Range("B22").Value = 1.044354791962178
$ws.Range("C22").Value = 0.1403306563405735
$ws.Range("D22").Value = 0.1240065706503799
$ws.Range("E22").Value = 0.1359734260365535
$ws.Range("F22").Value = 2.101697680817736
$ws.Range("J22").Value = 0.1858530138104619
$ws.Range("K22").Value = 0.614095883458333
$ws.Range("L22").Value = 0.2739056939461761
$ws.Range("O22").Value = 5.61200048669491

# Row 23
$ws.Range("B23").Value = 1.024827896290134
$ws.Range("C23").Value = 0.1399915565678569
$ws.Range("D23").Value = 0.1229833148976383
$ws.Range("E23").Value = 0.1358650788911682
$ws.Range("F23").Value = 2.103101975022952
$ws.Range("J23").Value = 0.186200084438024
$ws.Range("K23").Value = 0.5970068954898125
$ws.Range("L23").Value = 0.2716511385671225
$ws.Range("O23").Value = 5.619730401522958

# Row 24
$ws.Range("B24").Value = 0.9513408008294562
$ws.Range("C24").Value = 0.1387132299187854
$ws.Range("D24").Value = 0.1191733709785439
$ws.Range("E24").Value = 0.1355331489730318
$ws.Range("F24").Value = 2.110019549906667
$ws.Range("J24").Value = 0.1876183318356226
$ws.Range("K24").Value = 0.5324073182047755
$ws.Range("L24").Value = 0.2632722643516416
$ws.Range("O24").Value = 5.653376255291874

# Row 25
$ws.Range("B25").Value = 0.8731002884532586
$ws.Range("C25").Value = 0.1373478857392882
$ws.Range("D25").Value = 0.1152058520579402
$ws.Range("E25").Value = 0.1353436344023997
$ws.Range("F25").Value = 2.120902909898341
$ws.Range("J25").Value = 0.1893704493169786
$ws.Range("K25").Value = 0.4630086826598472
$ws.Range("L25").Value = 0.2545800857198373
$ws.Range("O25").Value = 5.699011697198898

Write-Host "Updated pl_mw data rows 2-25 for case with 380 kV"